$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 10 ("고재형") : fill in the previously-empty score cells ---
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 1
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0
$ws.Range("T10").Value = 0

# --- Row 11 ("김윤호") : fill in the previously-empty score cells ---
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 1
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 1
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = 1
$ws.Range("S11").Value = 1
$ws.Range("T11").Value = 1

# The SUM formulas in row 26 (C26:T26) automatically recalculate from the
# newly-entered values above, matching the totals recorded in the target file.

# --- Update the active selection left behind in the saved view ---
$ws.Range("I16").Select()

# --- Best-effort: restore the saved window geometry (host-environment dependent) ---
$excel.ActiveWindow.WindowState = -4143
$excel.Windows.Item(1).Left = 0
$excel.Windows.Item(1).Top = 0
$excel.Windows.Item(1).Width = 15630
$excel.Windows.Item(1).Height = 6945
